$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Azure SQL Database" bullet: merge the trailing " " + "is a cloud-based
#    relational database service" + "." runs into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" is a cloud-based relational database service.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " is a cloud-based relational database service.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Victor Gyokeres match-log bullets: drop the stray "Number #9" / "Sweden"
#    sample entries, keep the schema list (Opponent / Starts game? / Minutes
#    played). Delete the old duplicate "Opponent"/"Starts game?" paragraphs
#    first (indices 12 then 11, highest first so earlier indices stay valid),
#    then rename "Sweden" -> "Starts game?" and "Number #9" -> "Opponent".
# ---------------------------------------------------------------------------
$d.Paragraphs(12).Range.Delete() | Out-Null
$d.Paragraphs(11).Range.Delete() | Out-Null
$d.Paragraphs(10).Range.Text = "Starts game?"
$d.Paragraphs(9).Range.Text = "Opponent"

# ---------------------------------------------------------------------------
# 3) Re-pagination: the surrounding edits move where Word last rendered a
#    page break, so <w:lastRenderedPageBreak/> hops from the first run of
#    each pair to the next one. Do the text-preserving Find/Replace first
#    (this is what drops the stale <w:lastRenderedPageBreak/> and also fuses
#    any sibling runs that carried identical formatting back into one run),
#    then re-insert <w:lastRenderedPageBreak/> on the new host paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Successful Take-Ons", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Successful Take-Ons", 2) | Out-Null
$d.Content.Find.Execute("Passes Completed", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Passes Completed", 2) | Out-Null
$d.Content.Find.Execute("Carries into Final 3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Carries into Final 3", 2) | Out-Null

function Insert-LastRenderedPageBreak($paraRange, $pPrXml, $rPrXml, $text, $preserve) {
    $spacer = ""
    if ($preserve) { $spacer = ' xml:space="preserve"' }
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + `
        '<w:pPr>' + $pPrXml + '</w:pPr>' + `
        '<w:r>' + $rPrXml + '<w:lastRenderedPageBreak/><w:t' + $spacer + '>' + $text + '</w:t></w:r>' + `
        '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $paraRange.InsertXML($pkg) | Out-Null
}

$pPr_total = '<w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPr_24 = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Total`r") {
        Insert-LastRenderedPageBreak $p.Range $pPr_total $rPr_24 "Total" $false
        break
    }
}

$pPr_passesBlocked = '<w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Passes Blocked`r") {
        Insert-LastRenderedPageBreak $p.Range $pPr_passesBlocked $rPr_24 "Passes Blocked" $false
        break
    }
}

$pPr_miscontrols = '<w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="7"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Miscontrols `r") {
        Insert-LastRenderedPageBreak $p.Range $pPr_miscontrols $rPr_24 "Miscontrols " $true
        break
    }
}
